# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and "Correspond Handback
# DateTime" (col H) for row 3 (the b9d5d801... file) on both the zh-cn and
# de-de detail sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 08:38:56"
$wsZhCn.Range("H3").Value = "2016-03-19 08:39:16"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 08:38:59"
$wsDeDe.Range("H3").Value = "2016-03-19 08:39:21"
